{"js": "// The document contains a single 20-row x 5-column table of simple\n// arithmetic problems (e.g. \"27+45=\"). The edit replaces the text of\n// every cell, in row-major (reading) order, with a new problem. The\n// table keeps the same 20x5 shape before and after the edit \u2014 only the\n// cell contents change \u2014 so the whole change can be expressed as one\n// assignment to Table.values.\n\nconst newValues = [\n  [\"16+15=\", \"40-4=\", \"62-56=\", \"28+25=\", \"50-33=\"],\n  [\"65+26=\", \"40-5=\", \"28+14=\", \"17+18=\", \"52-44=\"],\n  [\"51-19=\", \"42-23=\", \"18+53=\", \"36-8=\", \"18+63=\"],\n  [\"75-17=\", \"75-18=\", \"49+12=\", \"60-42=\", \"62-44=\"],\n  [\"17+6=\", \"73-48=\", \"57-49=\", \"3+68=\", \"71-56=\"],\n  [\"81-9=\", \"90-41=\", \"18+49=\", \"64-28=\", \"14+18=\"],\n  [\"81-14=\", \"9+18=\", \"40-36=\", \"32-14=\", \"28+9=\"],\n  [\"56-28=\", \"87-9=\", \"98-29=\", \"82-74=\", \"8+57=\"],\n  [\"30-1=\", \"48+25=\", \"49+26=\", \"72-18=\", \"43-19=\"],\n  [\"58+5=\", \"15+47=\", \"44+19=\", \"43+49=\", \"96-29=\"],\n  [\"8+9=\", \"19+22=\", \"38+13=\", \"60-29=\", \"54-36=\"],\n  [\"35+19=\", \"73-67=\", \"77-29=\", \"52-15=\", \"9+87=\"],\n  [\"15+49=\", \"23+38=\", \"17+47=\", \"80-42=\", \"40-5=\"],\n  [\"24-5=\", \"21-3=\", \"45+47=\", \"75-16=\", \"61-25=\"],\n  [\"13+49=\", \"60-8=\", \"80-36=\", \"75+16=\", \"43-7=\"],\n  [\"45-38=\", \"23+28=\", \"83-24=\", \"85-68=\", \"36+38=\"],\n  [\"24-18=\", \"39+48=\", \"95-86=\", \"71-35=\", \"82-68=\"],\n  [\"76-37=\", \"8+64=\", \"44-36=\", \"72-26=\", \"16+25=\"],\n  [\"29+7=\", \"57-8=\", \"54+28=\", \"17-8=\", \"88-59=\"],\n  [\"62-37=\", \"95-87=\", \"62-29=\", \"79+13=\", \"82-27=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Sanity: bail out loudly (rather than silently mis-editing) if the\n// table shape ever doesn't match what we expect.\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    `Unexpected row count ${table.rowCount}, expected ${newValues.length}`\n  );\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table of simple\n# arithmetic problems (e.g. \"27+45=\"). The edit replaces the text of\n# every cell, in row-major (reading) order, with a new problem. The\n# table keeps the same 20x5 shape before and after the edit \u2014 only the\n# cell contents change \u2014 so the whole change can be expressed as writing\n# a new value into every Cell(row, col).\n\n$newValues = @(\n  @(\"16+15=\", \"40-4=\", \"62-56=\", \"28+25=\", \"50-33=\"),\n  @(\"65+26=\", \"40-5=\", \"28+14=\", \"17+18=\", \"52-44=\"),\n  @(\"51-19=\", \"42-23=\", \"18+53=\", \"36-8=\", \"18+63=\"),\n  @(\"75-17=\", \"75-18=\", \"49+12=\", \"60-42=\", \"62-44=\"),\n  @(\"17+6=\", \"73-48=\", \"57-49=\", \"3+68=\", \"71-56=\"),\n  @(\"81-9=\", \"90-41=\", \"18+49=\", \"64-28=\", \"14+18=\"),\n  @(\"81-14=\", \"9+18=\", \"40-36=\", \"32-14=\", \"28+9=\"),\n  @(\"56-28=\", \"87-9=\", \"98-29=\", \"82-74=\", \"8+57=\"),\n  @(\"30-1=\", \"48+25=\", \"49+26=\", \"72-18=\", \"43-19=\"),\n  @(\"58+5=\", \"15+47=\", \"44+19=\", \"43+49=\", \"96-29=\"),\n  @(\"8+9=\", \"19+22=\", \"38+13=\", \"60-29=\", \"54-36=\"),\n  @(\"35+19=\", \"73-67=\", \"77-29=\", \"52-15=\", \"9+87=\"),\n  @(\"15+49=\", \"23+38=\", \"17+47=\", \"80-42=\", \"40-5=\"),\n  @(\"24-5=\", \"21-3=\", \"45+47=\", \"75-16=\", \"61-25=\"),\n  @(\"13+49=\", \"60-8=\", \"80-36=\", \"75+16=\", \"43-7=\"),\n  @(\"45-38=\", \"23+28=\", \"83-24=\", \"85-68=\", \"36+38=\"),\n  @(\"24-18=\", \"39+48=\", \"95-86=\", \"71-35=\", \"82-68=\"),\n  @(\"76-37=\", \"8+64=\", \"44-36=\", \"72-26=\", \"16+25=\"),\n  @(\"29+7=\", \"57-8=\", \"54+28=\", \"17-8=\", \"88-59=\"),\n  @(\"62-37=\", \"95-87=\", \"62-29=\", \"79+13=\", \"82-27=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nif ($t.Rows.Count -ne $newValues.Length) {\n  throw \"Unexpected row count $($t.Rows.Count), expected $($newValues.Length)\"\n}\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$r - 1][$c - 1]\n  }\n}\n"}
